$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "268.79"

Set-TextValue $ws.Range("D3") "22.92"

Set-TextValue $ws.Range("D5") "0.06182"

Set-TextValue $ws.Range("D7") "6.669"

Set-TextValue $ws.Range("D8") "1.389"

Set-TextValue $ws.Range("D9") "0.8295"

Set-TextValue $ws.Range("D10") "0.01374"

Set-TextValue $ws.Range("D11") "0.1608"

Set-TextValue $ws.Range("D12") "0.08350"

Set-TextValue $ws.Range("D13") "0.03388"

Set-TextValue $ws.Range("D14") "0.03194"

Set-TextValue $ws.Range("D15") "0.09330"

Set-TextValue $ws.Range("D16") "3.840"

Set-TextValue $ws.Range("D17") "0.001637"

Set-TextValue $ws.Range("D18") "0.04727"

Set-TextValue $ws.Range("D19") "0.006398"

Set-TextValue $ws.Range("D20") "0.005666"

Set-TextValue $ws.Range("D21") "0.001077"

Set-TextValue $ws.Range("D22") "0.0001500"

Set-TextValue $ws.Range("D23") "3.726"

Set-TextValue $ws.Range("D24") "2.413"

Set-TextValue $ws.Range("D26") "0.1238"

Set-TextValue $ws.Range("D27") "0.0002703"

Set-TextValue $ws.Range("D40") "0.04702"

Set-TextValue $ws.Range("D41") "0.006946"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1161"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003301"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue $ws.Range("D44") "0.01192"

Set-TextValue $ws.Range("D45") "0.00006246"

Set-TextValue $ws.Range("D47") "0.00000000750"

Set-TextValue $ws.Range("D48") "0.9197"

$ws.Range("B49").Value = "CryptobidCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
Set-TextValue $ws.Range("D49") "0.00001400"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"

$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws.Range("D50") "0.002303"
$ws.Range("E50").Value = "49BOLOBOLO"

Set-TextValue $ws.Range("D51") "0.01240"
